$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-8) effectively rotate: each row now shows the data that
# previously belonged to the row two below it (wrapping around within rows 2-8).
# Apply the final target values directly, per the diff.

# Row 2
$ws.Range("D2").Value = 44424
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 75
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 18000
$ws.Range("N2").Value = "$/caja 15 kilos"
$ws.Range("P2").Value = 1200
$ws.Range("Q2").Value = 15

# Row 3
$ws.Range("D3").Value = 44424
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("N3").Value = "$/caja 15 kilos"
$ws.Range("P3").Value = 800
$ws.Range("Q3").Value = 15

# Row 4
$ws.Range("D4").Value = 44235
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("N4").Value = "$/bandeja 18 kilos"
$ws.Range("P4").Value = 778
$ws.Range("Q4").Value = 18

# Row 5
$ws.Range("D5").Value = 44235
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 70
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = "$/bandeja 18 kilos"
$ws.Range("P5").Value = 667
$ws.Range("Q5").Value = 18

# Row 6
$ws.Range("I6").Value = "Tercera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("P6").Value = 556

# Row 7
$ws.Range("D7").Value = 44238
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("P7").Value = 722

# Row 8
$ws.Range("D8").Value = 44238
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 11000
$ws.Range("P8").Value = 611
